$wb = $excel.ActiveWorkbook

# --- Keywords sheet: update row 7 keyword, add row 8 (new keyword/thesaurus pair) ---
$ws = $wb.Worksheets.Item("Keywords")
$ws.Range("A7").Value = "gross primary production"
$ws.Range("A8").Value = "Northeast U.S. Continental Shelf"
$ws.Range("B8").Value = "NOAA Large Marine Ecosystems"

# Make Keywords the active/selected sheet and select A7:B8, matching the saved view state
$ws.Activate()
$ws.Range("A7:B8").Select()
